# CU-Comprar Portátiles.docx
#
# 1) "El cliente rellena los datos y elige la opción Guardar." ->
#    "...elige la opción Comprar." split into 3 runs:
#       "El cliente rellena los datos y elige la opción " | "Comprar" | "."
#
# 2) "Flujo Alternativo 3 – al Paso " (+ existing run "4") ->
#    "Flujo Alternativo 2 – al Paso " (+ existing run "4") split into:
#       "Flujo Alternativo " | "2" | " – al Paso "
#
# 3) "Flujo Alternativo 4 – al Paso 5" -> "Flujo Alternativo 3 – al Paso 5"
#       "Flujo Alternativo " | "3" | " – al Paso 5"
#
# 4) "Flujo Alternativo 5 – al Paso 7" -> "Flujo Alternativo 4 – al Paso 7"
#       "Flujo Alternativo " | "4" | " – al Paso 7"
#
# 5) "Flujo Alternativo 6 – al Paso 7" -> "Flujo Alternativo 5 – al Paso 7"
#       "Flujo Alternativo " | "5" | " – al Paso 7"

$d = $word.ActiveDocument

# Forces a run boundary at the edges of $range by toggling (and reverting) a
# character property on it. Word COM does not expose "Run" objects directly
# (only Range), so a genuine format round-trip is the only reliable way to
# make the engine materialize a separate <w:r> for this span instead of
# silently re-merging it with its identically-formatted neighbour.
function Force-RunBoundary($range) {
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

function Split-ParagraphRuns($paragraph, [string[]]$segments) {
    # $segments must concatenate to exactly $paragraph.Range.Text (minus the
    # trailing paragraph mark). Re-materializes the paragraph's run layout so
    # each segment becomes its own <w:r>, preserving segment order/content.
    $pStart = $paragraph.Range.Start
    $offset = $pStart
    foreach ($seg in $segments) {
        $segStart = $offset
        $segEnd = $offset + $seg.Length
        if ($segStart -ne $pStart) {
            $segRange = $d.Range($segStart, $segEnd)
            Force-RunBoundary $segRange
        }
        $offset = $segEnd
    }
}

# --- 1) "... elige la opción Guardar." -> "... elige la opción Comprar." ---
$p8 = $d.Paragraphs(8)
$p8Start = $p8.Range.Start
$full8 = $p8.Range.Text
$oldTail8 = "Guardar."
$idx8 = $full8.IndexOf($oldTail8)
$tailStart8 = $p8Start + $idx8
$tailEnd8 = $tailStart8 + $oldTail8.Length
$d.Range($tailStart8, $tailEnd8).Text = "Comprar."

$prefix8 = "El cliente rellena los datos y elige la opción "
Split-ParagraphRuns $p8 @($prefix8, "Comprar", ".")

# --- 2) "Flujo Alternativo 3 – al Paso " (+ existing "4") -> "...2..." ---
$p16 = $d.Paragraphs(16)
$p16Start = $p16.Range.Start
$full16 = $p16.Range.Text
$oldNum16 = "3"
$idx16 = $full16.IndexOf($oldNum16)
$numStart16 = $p16Start + $idx16
$numEnd16 = $numStart16 + $oldNum16.Length
$d.Range($numStart16, $numEnd16).Text = "2"

$prefix16 = "Flujo Alternativo "
$mid16 = "2"
$suffix16 = " – al Paso "
Split-ParagraphRuns $p16 @($prefix16, $mid16, $suffix16)

# --- 3) "Flujo Alternativo 4 – al Paso 5" -> "...3..." ---
$p18 = $d.Paragraphs(18)
$p18Start = $p18.Range.Start
$full18 = $p18.Range.Text
$oldNum18 = "4"
$idx18 = $full18.IndexOf($oldNum18)
$numStart18 = $p18Start + $idx18
$numEnd18 = $numStart18 + $oldNum18.Length
$d.Range($numStart18, $numEnd18).Text = "3"
Split-ParagraphRuns $p18 @("Flujo Alternativo ", "3", " – al Paso 5")

# --- 4) "Flujo Alternativo 5 – al Paso 7" -> "...4..." ---
$p20 = $d.Paragraphs(20)
$p20Start = $p20.Range.Start
$full20 = $p20.Range.Text
$oldNum20 = "5"
$idx20 = $full20.IndexOf($oldNum20)
$numStart20 = $p20Start + $idx20
$numEnd20 = $numStart20 + $oldNum20.Length
$d.Range($numStart20, $numEnd20).Text = "4"
Split-ParagraphRuns $p20 @("Flujo Alternativo ", "4", " – al Paso 7")

# --- 5) "Flujo Alternativo 6 – al Paso 7" -> "...5..." ---
$p22 = $d.Paragraphs(22)
$p22Start = $p22.Range.Start
$full22 = $p22.Range.Text
$oldNum22 = "6"
$idx22 = $full22.IndexOf($oldNum22)
$numStart22 = $p22Start + $idx22
$numEnd22 = $numStart22 + $oldNum22.Length
$d.Range($numStart22, $numEnd22).Text = "5"
Split-ParagraphRuns $p22 @("Flujo Alternativo ", "5", " – al Paso 7")

Write-Host "Done"
